$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Paragraph 3 ("At this project, ... lose the race."):
# collapse the ten runs into one run and fix two small wording issues:
#   "not straight forward road"  -> "not a straight road"
#   "width of road is variable"  -> "width of the road is variable"
# --------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.MoveEnd(1, -1) | Out-Null
$p3.Text = "At this project, there are two robots and these robots try to catch each other. Then, it which catches other wins. Robots will compete on a road. This road is not a straight road. This road has elliptical shape and width of the road is variable. The aim of the project is that a robot should approach the other around 5 cm. While a robot tries to catch the other, it should not get out of line. If it is, it will lose the race. "

# --------------------------------------------------------------------
# Paragraph 4 ("Moreover, ... microcontroller"):
#   - drop the first-line indent
#   - rewrite/expand the text and re-split it into three runs:
#       "Moreover, ... used well. Also, ... should be used. Moreover, "
#       "to"
#       " finish the road, ... microcontroller."
# --------------------------------------------------------------------
$p4Para = $d.Paragraphs(4)
$p4Para.Format.FirstLineIndent = 0

$p4 = $p4Para.Range
$p4.MoveEnd(1, -1) | Out-Null

$p4run1 = "Moreover, they who control the motor properly will win since if they control the motor properly, the robot does not get out of line and robot goes faster. Because of this, I think PID controller should be used well. Also, to distinguish the road, image processing should be used. Moreover, "
$p4run2 = "to"
$p4run3 = " finish the road, while the robot is turning the upper part of the ellipse, it should slow down since if it is not, it can be driven away. Also, we will be familiar with some microcontroller."
$p4.Text = $p4run1 + $p4run2 + $p4run3

$p4Start = $d.Paragraphs(4).Range.Start
$p4Boundary1 = $p4Start + $p4run1.Length
$p4Boundary2 = $p4Boundary1 + $p4run2.Length

# Force a run break at each boundary by nudging a formatting property to a
# different value and straight back - this splits the run without altering
# the final appearance.
$p4Split = $d.Range($p4Boundary1, $p4Boundary2)
$p4Split.Font.Size = 11
$p4Split = $d.Range($p4Boundary1, $p4Boundary2)
$p4Split.Font.Size = 12

# --------------------------------------------------------------------
# Paragraph 5 ("<tab>Finally, ... peripheral.  "):
#   - keep the leading tab run untouched
#   - expand "Finally," into the full sentence
#   - move the _GoBack bookmark to sit right after that run
#   - the old trailing run becomes a single trailing space
# --------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5Start = $p5.Range.Start
$p5End = $p5.Range.End

$p5TextStart = $p5Start + 1   # skip the leading tab character
$p5TextEnd = $p5End - 1       # exclude the paragraph mark

$p5Finally = "Finally, we will learn a lot of things. These are image processing algorithm, PID controller, usage of the microcontroller with some peripheral.  "
$p5Trailing = " "

$p5Range = $d.Range($p5TextStart, $p5TextEnd)
$p5Range.Text = $p5Finally + $p5Trailing

$p5Boundary = $p5TextStart + $p5Finally.Length
$p5NewEnd = $p5Boundary + $p5Trailing.Length

$p5Split = $d.Range($p5Boundary, $p5NewEnd)
$p5Split.Font.Size = 11
$p5Split = $d.Range($p5Boundary, $p5NewEnd)
$p5Split.Font.Size = 12

# Re-create the _GoBack bookmark right between the two runs. Word only
# allows one bookmark per name, so this also removes it from its old
# location at the end of the document.
$bmRange = $d.Range($p5Boundary, $p5Boundary)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
